$d = $word.ActiveDocument

# Standard run properties used throughout this document's list items.
$rpr = '<w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US" w:eastAsia="vi-VN"/></w:rPr>'

$pkgHeader = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>'
$pkgFooter = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

function Replace-ParagraphTail {
    # Finds $anchorText -- a globally-unique search string whose match
    # ends exactly where its paragraph ends -- then replaces only its
    # last $tailLen characters of that match with the raw run-level OOXML
    # in $innerXml (wrapped in a throwaway <w:p> so InsertXML accepts it).
    # $tailLen must always be passed explicitly (no relying on defaults).
    param(
        [string]$anchorText,
        [string]$innerXml,
        [int]$tailLen
    )
    $f = $d.Content
    $found = $f.Find.Execute($anchorText, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        Write-Host "NOT FOUND: $anchorText"
        return
    }
    $tailStart = $f.End - $tailLen
    $target = $d.Range($tailStart, $f.End)
    $xml = $pkgHeader + '<w:p>' + $innerXml + '</w:p>' + $pkgFooter
    $target.InsertXML($xml)
}

# 1) "Đặc tả yêu cầu hệ thống sẽ thực hiện: Vũ, Tài, N.Phúc"
#    -> "...Vũ, Tài, " + gramStart + "N.Phúc" + gramEnd
Replace-ParagraphTail "Vũ, Tài, N.Phúc" `
    ('<w:proofErr w:type="gramStart"/><w:r>' + $rpr + '<w:t>N.Phúc</w:t></w:r><w:proofErr w:type="gramEnd"/>') `
    6

# 2) "Công cụ: Power Designer: V.Phúc, Lâm, Anh"
#    -> "...Power Designer: " + gramStart + "V.Phúc" + gramEnd + ", Lâm, Anh"
Replace-ParagraphTail "V.Phúc, Lâm, Anh" `
    ('<w:proofErr w:type="gramStart"/><w:r>' + $rpr + '<w:t>V.Phúc</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r>' + $rpr + '<w:t>, Lâm, Anh</w:t></w:r>') `
    16

# 3) "Thành, Thanh, N.Phúc" -> "Thành, Thanh, " + gramStart + "N.Phúc" + gramEnd
Replace-ParagraphTail "Thanh, N.Phúc" `
    ('<w:proofErr w:type="gramStart"/><w:r>' + $rpr + '<w:t>N.Phúc</w:t></w:r><w:proofErr w:type="gramEnd"/>') `
    6

# 4) "Chỉnh sửa word: Lâm, V.Phúc" (already split into "V." + "Phúc" runs)
#    -> wrap those two existing runs with gramStart / gramEnd, unchanged text
Replace-ParagraphTail "Lâm, V.Phúc" `
    ('<w:proofErr w:type="gramStart"/><w:r>' + $rpr + '<w:t>V.</w:t></w:r><w:r>' + $rpr + '<w:t>Phúc</w:t></w:r><w:proofErr w:type="gramEnd"/>') `
    6

# 5) "Phần 4: N.Phúc, Anh" -> "Phần 4: " + gramStart + "N.Phúc" + gramEnd + ", Anh"
Replace-ParagraphTail "Phần 4: N.Phúc, Anh" `
    ('<w:proofErr w:type="gramStart"/><w:r>' + $rpr + '<w:t>N.Phúc</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r>' + $rpr + '<w:t>, Anh</w:t></w:r>') `
    19

# 6) "Trả lời các câu hỏi và battlerap: Phát, Phúc, Thanh"
#    -> "Trả lời các câu hỏi và " + "phản biện" + ": Phát, Phúc, Thanh"  (no proofErr marks)
Replace-ParagraphTail "battlerap: Phát, Phúc, Thanh" `
    ('<w:r>' + $rpr + '<w:t>phản biện</w:t></w:r><w:r>' + $rpr + '<w:t>: Phát, Phúc, Thanh</w:t></w:r>') `
    28

# 7) Insert a brand-new list paragraph right after the
#    "Công cụ thực hiện: VS code, VS studio, SQL server, Github" item.
$pPr = '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="46"/></w:numPr><w:spacing w:after="60" w:line="240" w:lineRule="auto"/><w:jc w:val="left"/>' + $rpr + '</w:pPr>'

foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text -like "*Công cụ thực hiện: VS code, VS studio, SQL server, Github*") {
        $target = $d.Range($para.Range.Start, $para.Range.End)
        $xml = $pkgHeader `
            + '<w:p>' + $pPr + '<w:r>' + $rpr + '<w:t>Công cụ thực hiện: VS code, VS studio, SQL server, Github</w:t></w:r></w:p>' `
            + '<w:p>' + $pPr + '<w:r>' + $rpr + '<w:t>Thực hiện code các chức năng: Cả nhóm</w:t></w:r></w:p>' `
            + $pkgFooter
        $target.InsertXML($xml)
        break
    }
}
